# DEV-5416: enable references to other ontologies in resource cardinalities
#
# 1) "Owner" sheet: append a cardinality row that references a property
#    from another ontology ("other-onto:hasProp"), cardinality "0-1",
#    gui_order 14.
# 2) "Image" sheet: append two new cardinality rows ("seqnum" and
#    "isPartOf", both cardinality "0-1") that are needed so other
#    resources can reference Image via isPartOf/seqnum.
# 3) Leave "Image" as the active/selected sheet+cell, which also clears
#    the previous tab-selection on "classes".

$wb = $excel.ActiveWorkbook

# --- Owner sheet: new cardinality row 17 -----------------------------
$wsOwner = $wb.Worksheets.Item("Owner")
$wsOwner.Range("A17").Value = "other-onto:hasProp"
$wsOwner.Range("B17").Value = "0-1"
$wsOwner.Range("C17").Value = 14
$wsOwner.Range("A17").Select() | Out-Null

# --- Image sheet: new cardinality rows 3 and 4 ------------------------
$wsImage = $wb.Worksheets.Item("Image")
$wsImage.Range("A3").Value = "seqnum"
$wsImage.Range("B3").Value = "0-1"
$wsImage.Range("A4").Value = "isPartOf"
$wsImage.Range("B4").Value = "0-1"

# Make "Image" the active tab/sheet, with A5 selected, matching the
# saved workbook state after the edit.
$wsImage.Activate() | Out-Null
$wsImage.Range("A5").Select() | Out-Null
